$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# EPBDS-8247: add two new "Collect as Double and Double Map" rule tables
# (rows 37-42 and rows 48-53), each mirroring the layout already used for
# the existing tables (merged, centered header row + a bordered data grid).
# ---------------------------------------------------------------------------

# --- Block 1: rows 37-42 (SimpleRules ... testSimpleRulesMap2) -------------
$ws.Range("B6").Copy()
$ws.Range("B37:E37").PasteSpecial(-4122)
$ws.Range("B37").Value = "SimpleRules  Collect as Double and Double Map testSimpleRulesMap2 (Integer a, String b)"

$ws.Range("B7").Copy()
$ws.Range("B38:E38").PasteSpecial(-4122)
$ws.Range("B38").Value = "A"
$ws.Range("C38").Value = "B"
$ws.Range("D38").Value = "A"
$ws.Range("E38").Value = "B"

$ws.Range("B8").Copy()
$ws.Range("B39:E39").PasteSpecial(-4122)
$ws.Range("B39").Value = 1
$ws.Range("C39").Value = "test1"
$ws.Range("D39").Value = "1.2"
$ws.Range("E39").Value = "1.2"

$ws.Range("B8").Copy()
$ws.Range("B40:E40").PasteSpecial(-4122)
$ws.Range("B40").Value = 1
$ws.Range("C40").Value = "test1"
$ws.Range("D40").Value = "1.22"
$ws.Range("E40").Value = "1.22"

$ws.Range("B8").Copy()
$ws.Range("B41:E41").PasteSpecial(-4122)
$ws.Range("B41").Value = 2
$ws.Range("C41").Value = "test2"
$ws.Range("D41").Value = "1.3"
$ws.Range("E41").Value = "1.3"

$ws.Range("B8").Copy()
$ws.Range("B42:E42").PasteSpecial(-4122)
$ws.Range("B42").Value = 1
$ws.Range("C42").Value = "test1"
$ws.Range("D42").Value = "1.1"
$ws.Range("E42").Value = "1.1"

# --- Block 2: rows 48-53 (SmartRules ... testSimpleRulesMap1) --------------
$ws.Range("B6").Copy()
$ws.Range("B48:E48").PasteSpecial(-4122)
$ws.Range("B48").Value = "SmartRules  Collect as Double and Double Map testSimpleRulesMap1 (Integer a, String b)"

$ws.Range("B7").Copy()
$ws.Range("B49:E49").PasteSpecial(-4122)
$ws.Range("B49").Value = "A"
$ws.Range("C49").Value = "B"
$ws.Range("D49").Value = "A"
$ws.Range("E49").Value = "B"

$ws.Range("B8").Copy()
$ws.Range("B50:E50").PasteSpecial(-4122)
$ws.Range("B50").Value = 1
$ws.Range("C50").Value = "test1"
$ws.Range("D50").Value = "1.2"
$ws.Range("E50").Value = "1.2"

$ws.Range("B8").Copy()
$ws.Range("B51:E51").PasteSpecial(-4122)
$ws.Range("B51").Value = 1
$ws.Range("C51").Value = "test1"
$ws.Range("D51").Value = "1.22"
$ws.Range("E51").Value = "1.22"

$ws.Range("B8").Copy()
$ws.Range("B52:E52").PasteSpecial(-4122)
$ws.Range("B52").Value = 2
$ws.Range("C52").Value = "test2"
$ws.Range("D52").Value = "1.3"
$ws.Range("E52").Value = "1.3"

$ws.Range("B8").Copy()
$ws.Range("B53:E53").PasteSpecial(-4122)
$ws.Range("B53").Value = 1
$ws.Range("C53").Value = "test1"
$ws.Range("D53").Value = "1.1"
$ws.Range("E53").Value = "1.1"

# --- Merge the two new header rows, keeping a uniform centered/bordered ----
# --- style across the whole merged range (re-paste the header format ------
# --- after merging, since Merge() redistributes borders like real Excel). --
$ws.Range("B37:E37").Merge()
$ws.Range("B6").Copy()
$ws.Range("B37").PasteSpecial(-4122)
$ws.Range("C37").PasteSpecial(-4122)
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("E37").PasteSpecial(-4122)

$ws.Range("B48:E48").Merge()
$ws.Range("B6").Copy()
$ws.Range("B48").PasteSpecial(-4122)
$ws.Range("C48").PasteSpecial(-4122)
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("E48").PasteSpecial(-4122)

# --- Final view state: scrolled down, selection on I44 ---------------------
$ws.Range("I44").Select()
